$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 29 metrics (C,D,E,F) per updated data for "add e bibi"
$ws.Range("C29").Value = 206
$ws.Range("D29").Value = 33
$ws.Range("E29").Value = 173
$ws.Range("F29").Value = 5.679862306368331
